$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 12 de Octubre de 2020 a las 10:26"

# Row 7 - Rusia
$ws.Range("B7").Value = 1312310
$ws.Range("C7").Value = 13592
$ws.Range("D7").Value = 1024235
$ws.Range("E7").Value = 265353
$ws.Range("G7").Value = 125
$ws.Range("H7").Value = 22722

# Row 21 - Filipinas
$ws.Range("B21").Value = 342816
$ws.Range("C21").Value = 3564
$ws.Range("D21").Value = 293152
$ws.Range("E21").Value = 43332
$ws.Range("G21").Value = 11
$ws.Range("H21").Value = 6332

# Row 37 - Polonia
$ws.Range("D37").Value = 81201
$ws.Range("E37").Value = 41611

# Row 62 - Singapur
$ws.Range("B62").Value = 57880
$ws.Range("C62").Value = 4
$ws.Range("E62").Value = 148

# Row 90 - Croacia
$ws.Range("B90").Value = 20621
$ws.Range("C90").Value = 181
$ws.Range("D90").Value = 17582
$ws.Range("E90").Value = 2712
$ws.Range("G90").Value = 3
$ws.Range("H90").Value = 327

# Row 140 - Estonia
$ws.Range("B140").Value = 3883
$ws.Range("C140").Value = 18
$ws.Range("D140").Value = 2967
$ws.Range("E140").Value = 848

# Row 207 - Timor Oriental
$ws.Range("B207").Value = 29
$ws.Range("C207").Value = 1
$ws.Range("E207").Value = 1
